$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (word dropped from the report: "toy")
$ws.Rows.Item(40).Delete()

# Row 1
$ws.Cells.Item(1, 1).Value = "negative"
$ws.Cells.Item(1, 10).Value = "positive"

# Row 2
$ws.Cells.Item(2, 1).Value = "name"
$ws.Cells.Item(2, 2).Value = "anchor score"
$ws.Cells.Item(2, 3).Value = "type occurences"
$ws.Cells.Item(2, 4).Value = "total occurences"
$ws.Cells.Item(2, 5).Value = "+%"
$ws.Cells.Item(2, 6).Value = "-%"
$ws.Cells.Item(2, 7).Value = "both"
$ws.Cells.Item(2, 8).Value = "normal"
$ws.Cells.Item(2, 10).Value = "name"
$ws.Cells.Item(2, 11).Value = "anchor score"
$ws.Cells.Item(2, 12).Value = "type occurences"
$ws.Cells.Item(2, 13).Value = "total occurences"
$ws.Cells.Item(2, 14).Value = "+%"
$ws.Cells.Item(2, 15).Value = "-%"
$ws.Cells.Item(2, 16).Value = "both"
$ws.Cells.Item(2, 17).Value = "normal"

# Row 3
$ws.Cells.Item(3, 1).Value = "poorly"
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 46
$ws.Cells.Item(3, 4).Value = 46
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = "awesome"
$ws.Cells.Item(3, 11).Value = 0.8153846153846154
$ws.Cells.Item(3, 12).Value = 53
$ws.Cells.Item(3, 13).Value = 53
$ws.Cells.Item(3, 14).Value = 1
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = $false
$ws.Cells.Item(3, 17).Value = 12

# Row 4
$ws.Cells.Item(4, 1).Value = "disappointing"
$ws.Cells.Item(4, 2).Value = 0.8863636363636364
$ws.Cells.Item(4, 3).Value = 39
$ws.Cells.Item(4, 4).Value = 39
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Value = 5
$ws.Cells.Item(4, 10).Value = "wonderful"
$ws.Cells.Item(4, 11).Value = 0.8035714285714286
$ws.Cells.Item(4, 12).Value = 45
$ws.Cells.Item(4, 13).Value = 45
$ws.Cells.Item(4, 14).Value = 1
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = $false
$ws.Cells.Item(4, 17).Value = 11

# Row 5
$ws.Cells.Item(5, 1).Value = "broke"
$ws.Cells.Item(5, 2).Value = 0.7621359223300971
$ws.Cells.Item(5, 3).Value = 157
$ws.Cells.Item(5, 4).Value = 157
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = $false
$ws.Cells.Item(5, 8).Value = 49
$ws.Cells.Item(5, 10).Value = "favorite"
$ws.Cells.Item(5, 11).Value = 0.6989247311827957
$ws.Cells.Item(5, 12).Value = 65
$ws.Cells.Item(5, 13).Value = 65
$ws.Cells.Item(5, 14).Value = 1
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = $false
$ws.Cells.Item(5, 17).Value = 28

# Row 6
$ws.Cells.Item(6, 1).Value = "disappointed"
$ws.Cells.Item(6, 2).Value = 0.7526881720430108
$ws.Cells.Item(6, 3).Value = 140
$ws.Cells.Item(6, 4).Value = 140
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = 46
$ws.Cells.Item(6, 10).Value = "classic"
$ws.Cells.Item(6, 11).Value = 0.5471698113207547
$ws.Cells.Item(6, 12).Value = 29
$ws.Cells.Item(6, 13).Value = 29
$ws.Cells.Item(6, 14).Value = 1
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = $false
$ws.Cells.Item(6, 17).Value = 24

# Row 7
$ws.Cells.Item(7, 1).Value = "however"
$ws.Cells.Item(7, 2).Value = 0.75
$ws.Cells.Item(7, 3).Value = 48
$ws.Cells.Item(7, 4).Value = 48
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = 16
$ws.Cells.Item(7, 10).Value = "thank"
$ws.Cells.Item(7, 11).Value = 0.4927536231884058
$ws.Cells.Item(7, 12).Value = 34
$ws.Cells.Item(7, 13).Value = 34
$ws.Cells.Item(7, 14).Value = 1
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = $false
$ws.Cells.Item(7, 17).Value = 35

# Row 8
$ws.Cells.Item(8, 1).Value = "poor"
$ws.Cells.Item(8, 2).Value = 0.6901408450704225
$ws.Cells.Item(8, 3).Value = 49
$ws.Cells.Item(8, 4).Value = 49
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = 22
$ws.Cells.Item(8, 10).Value = "excellent"
$ws.Cells.Item(8, 11).Value = 0.46875
$ws.Cells.Item(8, 12).Value = 30
$ws.Cells.Item(8, 13).Value = 30
$ws.Cells.Item(8, 14).Value = 1
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = $false
$ws.Cells.Item(8, 17).Value = 34

# Row 9
$ws.Cells.Item(9, 1).Value = "waste"
$ws.Cells.Item(9, 2).Value = 0.6621621621621622
$ws.Cells.Item(9, 3).Value = 98
$ws.Cells.Item(9, 4).Value = 98
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = $false
$ws.Cells.Item(9, 8).Value = 50
$ws.Cells.Item(9, 10).Value = "great"
$ws.Cells.Item(9, 11).Value = 0.369672131147541
$ws.Cells.Item(9, 12).Value = 451
$ws.Cells.Item(9, 13).Value = 451
$ws.Cells.Item(9, 14).Value = 1
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = $false
$ws.Cells.Item(9, 17).Value = 769

# Row 10
$ws.Cells.Item(10, 1).Value = "junk"
$ws.Cells.Item(10, 2).Value = 0.6
$ws.Cells.Item(10, 3).Value = 33
$ws.Cells.Item(10, 4).Value = 33
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = 22
$ws.Cells.Item(10, 10).Value = "love"
$ws.Cells.Item(10, 11).Value = 0.3314203730272597
$ws.Cells.Item(10, 12).Value = 231
$ws.Cells.Item(10, 13).Value = 231
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = $false
$ws.Cells.Item(10, 17).Value = 466

# Row 11
$ws.Cells.Item(11, 1).Value = "smaller"
$ws.Cells.Item(11, 2).Value = 0.5966386554621849
$ws.Cells.Item(11, 3).Value = 71
$ws.Cells.Item(11, 4).Value = 71
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = $false
$ws.Cells.Item(11, 8).Value = 48
$ws.Cells.Item(11, 10).Value = "loves"
$ws.Cells.Item(11, 11).Value = 0.3091286307053942
$ws.Cells.Item(11, 12).Value = 149
$ws.Cells.Item(11, 13).Value = 149
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = $false
$ws.Cells.Item(11, 17).Value = 333

# Row 12
$ws.Cells.Item(12, 1).Value = "small"
$ws.Cells.Item(12, 2).Value = 0.5130434782608696
$ws.Cells.Item(12, 3).Value = 177
$ws.Cells.Item(12, 4).Value = 177
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = $false
$ws.Cells.Item(12, 8).Value = 168
$ws.Cells.Item(12, 10).Value = "best"
$ws.Cells.Item(12, 11).Value = 0.2583333333333334
$ws.Cells.Item(12, 12).Value = 31
$ws.Cells.Item(12, 13).Value = 31
$ws.Cells.Item(12, 14).Value = 1
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = $false
$ws.Cells.Item(12, 17).Value = 89

# Row 13
$ws.Cells.Item(13, 1).Value = "paint"
$ws.Cells.Item(13, 2).Value = 0.4603174603174603
$ws.Cells.Item(13, 3).Value = 29
$ws.Cells.Item(13, 4).Value = 29
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = $false
$ws.Cells.Item(13, 8).Value = 34
$ws.Cells.Item(13, 10).Value = "perfect"
$ws.Cells.Item(13, 11).Value = 0.2530120481927711
$ws.Cells.Item(13, 12).Value = 42
$ws.Cells.Item(13, 13).Value = 42
$ws.Cells.Item(13, 14).Value = 1
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = $false
$ws.Cells.Item(13, 17).Value = 124

# Row 14
$ws.Cells.Item(14, 1).Value = "plastic"
$ws.Cells.Item(14, 2).Value = 0.4488188976377953
$ws.Cells.Item(14, 3).Value = 57
$ws.Cells.Item(14, 4).Value = 57
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = $false
$ws.Cells.Item(14, 8).Value = 70
$ws.Cells.Item(14, 10).Value = "loved"
$ws.Cells.Item(14, 11).Value = 0.2415902140672783
$ws.Cells.Item(14, 12).Value = 79
$ws.Cells.Item(14, 13).Value = 79
$ws.Cells.Item(14, 14).Value = 1
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = $false
$ws.Cells.Item(14, 17).Value = 248

# Row 15
$ws.Cells.Item(15, 1).Value = "broken"
$ws.Cells.Item(15, 2).Value = 0.4457831325301205
$ws.Cells.Item(15, 3).Value = 37
$ws.Cells.Item(15, 4).Value = 37
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = $false
$ws.Cells.Item(15, 8).Value = 46
$ws.Cells.Item(15, 10).Value = "friends"
$ws.Cells.Item(15, 11).Value = 0.1904761904761905
$ws.Cells.Item(15, 12).Value = 36
$ws.Cells.Item(15, 13).Value = 36
$ws.Cells.Item(15, 14).Value = 1
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = $false
$ws.Cells.Item(15, 17).Value = 153

# Row 16
$ws.Cells.Item(16, 1).Value = "apart"
$ws.Cells.Item(16, 2).Value = 0.4210526315789473
$ws.Cells.Item(16, 3).Value = 40
$ws.Cells.Item(16, 4).Value = 40
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = $false
$ws.Cells.Item(16, 8).Value = 55
$ws.Cells.Item(16, 10).Value = "christmas"
$ws.Cells.Item(16, 11).Value = 0.1285140562248996
$ws.Cells.Item(16, 12).Value = 32
$ws.Cells.Item(16, 13).Value = 32
$ws.Cells.Item(16, 14).Value = 1
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = $false
$ws.Cells.Item(16, 17).Value = 217

# Row 17
$ws.Cells.Item(17, 1).Value = "ok"
$ws.Cells.Item(17, 2).Value = 0.3671875
$ws.Cells.Item(17, 3).Value = 47
$ws.Cells.Item(17, 4).Value = 47
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = $false
$ws.Cells.Item(17, 8).Value = 81
$ws.Cells.Item(17, 10).Value = "fun"
$ws.Cells.Item(17, 11).Value = 0.09569798068481124
$ws.Cells.Item(17, 12).Value = 109
$ws.Cells.Item(17, 13).Value = 111
$ws.Cells.Item(17, 14).Value = 0.98
$ws.Cells.Item(17, 15).Value = 0.02000000000000002
$ws.Cells.Item(17, 16).Value = $true
$ws.Cells.Item(17, 17).Value = 1030

# Row 18
$ws.Cells.Item(18, 1).Value = "difficult"
$ws.Cells.Item(18, 2).Value = 0.3483146067415731
$ws.Cells.Item(18, 3).Value = 31
$ws.Cells.Item(18, 4).Value = 31
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = $false
$ws.Cells.Item(18, 8).Value = 58
$ws.Cells.Item(18, 10).Value = "game"
$ws.Cells.Item(18, 11).Value = 0.05974025974025974
$ws.Cells.Item(18, 12).Value = 92
$ws.Cells.Item(18, 13).Value = 93
$ws.Cells.Item(18, 14).Value = 0.99
$ws.Cells.Item(18, 15).Value = 0.01000000000000001
$ws.Cells.Item(18, 16).Value = $true
$ws.Cells.Item(18, 17).Value = 1448

# Row 19
$ws.Cells.Item(19, 1).Value = "thought"
$ws.Cells.Item(19, 2).Value = 0.3316831683168317
$ws.Cells.Item(19, 3).Value = 67
$ws.Cells.Item(19, 4).Value = 67
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = $false
$ws.Cells.Item(19, 8).Value = 135

# Row 20
$ws.Cells.Item(20, 1).Value = "cheap"
$ws.Cells.Item(20, 2).Value = 0.3175355450236967
$ws.Cells.Item(20, 3).Value = 67
$ws.Cells.Item(20, 4).Value = 67
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = $false
$ws.Cells.Item(20, 8).Value = 144

# Row 21
$ws.Cells.Item(21, 1).Value = "though"
$ws.Cells.Item(21, 2).Value = 0.264957264957265
$ws.Cells.Item(21, 3).Value = 31
$ws.Cells.Item(21, 4).Value = 31
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = $false
$ws.Cells.Item(21, 8).Value = 86

# Row 22
$ws.Cells.Item(22, 1).Value = "size"
$ws.Cells.Item(22, 2).Value = 0.2525773195876289
$ws.Cells.Item(22, 3).Value = 49
$ws.Cells.Item(22, 4).Value = 49
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = $false
$ws.Cells.Item(22, 8).Value = 145

# Row 23
$ws.Cells.Item(23, 1).Value = "would"
$ws.Cells.Item(23, 2).Value = 0.2154531946508172
$ws.Cells.Item(23, 3).Value = 145
$ws.Cells.Item(23, 4).Value = 146
$ws.Cells.Item(23, 5).Value = 0.01
$ws.Cells.Item(23, 6).Value = 0.99
$ws.Cells.Item(23, 7).Value = $true
$ws.Cells.Item(23, 8).Value = 528

# Row 24
$ws.Cells.Item(24, 1).Value = "hard"
$ws.Cells.Item(24, 2).Value = 0.21
$ws.Cells.Item(24, 3).Value = 42
$ws.Cells.Item(24, 4).Value = 42
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = $false
$ws.Cells.Item(24, 8).Value = 158

# Row 25
$ws.Cells.Item(25, 1).Value = "item"
$ws.Cells.Item(25, 2).Value = 0.2028985507246377
$ws.Cells.Item(25, 3).Value = 56
$ws.Cells.Item(25, 4).Value = 56
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = $false
$ws.Cells.Item(25, 8).Value = 220

# Row 26
$ws.Cells.Item(26, 1).Value = "work"
$ws.Cells.Item(26, 2).Value = 0.1930379746835443
$ws.Cells.Item(26, 3).Value = 61
$ws.Cells.Item(26, 4).Value = 61
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = $false
$ws.Cells.Item(26, 8).Value = 255

# Row 27
$ws.Cells.Item(27, 1).Value = "money"
$ws.Cells.Item(27, 2).Value = 0.189873417721519
$ws.Cells.Item(27, 3).Value = 60
$ws.Cells.Item(27, 4).Value = 60
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = $false
$ws.Cells.Item(27, 8).Value = 256

# Row 28
$ws.Cells.Item(28, 1).Value = "could"
$ws.Cells.Item(28, 2).Value = 0.1847133757961783
$ws.Cells.Item(28, 3).Value = 29
$ws.Cells.Item(28, 4).Value = 29
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = $false
$ws.Cells.Item(28, 8).Value = 128

# Row 29
$ws.Cells.Item(29, 1).Value = "product"
$ws.Cells.Item(29, 2).Value = 0.1828193832599119
$ws.Cells.Item(29, 3).Value = 83
$ws.Cells.Item(29, 4).Value = 83
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = $false
$ws.Cells.Item(29, 8).Value = 371

# Row 30
$ws.Cells.Item(30, 1).Value = "used"
$ws.Cells.Item(30, 2).Value = 0.1771428571428571
$ws.Cells.Item(30, 3).Value = 31
$ws.Cells.Item(30, 4).Value = 31
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = $false
$ws.Cells.Item(30, 8).Value = 144

# Row 31
$ws.Cells.Item(31, 1).Value = "better"
$ws.Cells.Item(31, 2).Value = 0.1635514018691589
$ws.Cells.Item(31, 3).Value = 35
$ws.Cells.Item(31, 4).Value = 35
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = $false
$ws.Cells.Item(31, 8).Value = 179

# Row 32
$ws.Cells.Item(32, 1).Value = "price"
$ws.Cells.Item(32, 2).Value = 0.1408045977011494
$ws.Cells.Item(32, 3).Value = 49
$ws.Cells.Item(32, 4).Value = 49
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = $false
$ws.Cells.Item(32, 8).Value = 299

# Row 33
$ws.Cells.Item(33, 1).Value = "3"
$ws.Cells.Item(33, 2).Value = 0.1169354838709677
$ws.Cells.Item(33, 3).Value = 29
$ws.Cells.Item(33, 4).Value = 29
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = $false
$ws.Cells.Item(33, 8).Value = 219

# Row 34
$ws.Cells.Item(34, 1).Value = "2"
$ws.Cells.Item(34, 2).Value = 0.1161048689138577
$ws.Cells.Item(34, 3).Value = 31
$ws.Cells.Item(34, 4).Value = 31
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = $false
$ws.Cells.Item(34, 8).Value = 236

# Row 35
$ws.Cells.Item(35, 1).Value = "use"
$ws.Cells.Item(35, 2).Value = 0.1150684931506849
$ws.Cells.Item(35, 3).Value = 42
$ws.Cells.Item(35, 4).Value = 42
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = $false
$ws.Cells.Item(35, 8).Value = 323

# Row 36
$ws.Cells.Item(36, 1).Value = "like"
$ws.Cells.Item(36, 2).Value = 0.08881578947368421
$ws.Cells.Item(36, 3).Value = 54
$ws.Cells.Item(36, 4).Value = 54
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = $false
$ws.Cells.Item(36, 8).Value = 554

# Row 37
$ws.Cells.Item(37, 1).Value = "little"
$ws.Cells.Item(37, 2).Value = 0.08482142857142858
$ws.Cells.Item(37, 3).Value = 38
$ws.Cells.Item(37, 4).Value = 39
$ws.Cells.Item(37, 5).Value = 0.03
$ws.Cells.Item(37, 6).Value = 0.97
$ws.Cells.Item(37, 7).Value = $true
$ws.Cells.Item(37, 8).Value = 410

# Row 38
$ws.Cells.Item(38, 1).Value = "much"
$ws.Cells.Item(38, 2).Value = 0.07746478873239436
$ws.Cells.Item(38, 3).Value = 33
$ws.Cells.Item(38, 4).Value = 41
$ws.Cells.Item(38, 5).Value = 0.2
$ws.Cells.Item(38, 6).Value = 0.8
$ws.Cells.Item(38, 7).Value = $true
$ws.Cells.Item(38, 8).Value = 393

# Row 39
$ws.Cells.Item(39, 1).Value = "one"
$ws.Cells.Item(39, 2).Value = 0.04701397712833545
$ws.Cells.Item(39, 3).Value = 37
$ws.Cells.Item(39, 4).Value = 44
$ws.Cells.Item(39, 5).Value = 0.16
$ws.Cells.Item(39, 6).Value = 0.84
$ws.Cells.Item(39, 7).Value = $true
$ws.Cells.Item(39, 8).Value = 750
